$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.156.00'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.828.29'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('D4').Value = "'0.9983"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'242.59"
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = "'0.6198"
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').Value = "'0.9999"
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '1.828.80'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.07345"
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = "'0.2910"
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = "'23.24"
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.07626"
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.832.36'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'4.961"
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = "'0.6682"
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = "'82.32"
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.000008985"
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = "'5.834"
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '29.153.23'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.087.38'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = "'235.76"
$ws.Range('E21').Value = '  +1.70%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = "'12.45"
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = "'0.9996"
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = "'7.343"
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = "'1.000"
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'158.26"
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.1391"
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = "'8.519"
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'17.60"
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.486"
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.05879"
$ws.Range('E31').Value = '  +6.75%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'4.077"
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'4.081"
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').Value = "'1.209"
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.860"
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'0.7245"
$ws.Range('E37').Value = '  -1.72%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = "'2.609"
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.844"
$ws.Range('E39').Value = '  +2.59%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.225.87'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = "'0.01756"
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'6.219"
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'0.9075"
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'1.000"
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '1.991.47'
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = "'101.67"
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'65.55"
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'0.5035"
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('D49').Value = "'0.4026"
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = "'0.00000000117"
$ws.Range('E50').Value = '  -3.10%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'9.112"
$ws.Range('E51').Value = '  -0.24%  '
